$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Update the "Förändrad" (changed) date column C for rows 2-6
# from 2023-09-01 (45170) to 2023-09-05 (45174), keeping existing
# date formatting (style) intact.
for ($row = 2; $row -le 6; $row++) {
    $ws.Cells.Item($row, 3).Value = "2023-09-05"
}
